$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "job description" column (J) header + data
$ws.Range("J5").Value = "job description"
$ws.Range("J6").Value = "a"
$ws.Range("J7").Value = "s"
$ws.Range("J8").Value = "ds"
$ws.Range("J9").Value = "wd"
$ws.Range("J10").Value = "dsada"
$ws.Range("J11").Value = "dsada"
$ws.Range("J12").Value = "ds"

# Match the author's final selection / scroll position
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("J13").Select()
